$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 1.4
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 2.2
$ws.Range("O2").Value = 1.67
$ws.Range("P2").Value = 1.5
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.91
$ws.Range("S2").Value = 1.91
$ws.Range("X2").Value = 23
$ws.Range("Z2").Value = 8
$ws.Range("AA2").Value = 5.5
$ws.Range("AB2").Value = 15
$ws.Range("AD2").Value = 301
# Row 3
$ws.Range("L3").Value = 1.36
$ws.Range("M3").Value = 3.2
$ws.Range("Z3").Value = 8.5
# Row 4
$ws.Range("N4").Value = 1.85
$ws.Range("O4").Value = 2.05
# Row 5
$ws.Range("K5").Value = 5
# Row 6
$ws.Range("N6").Value = 2.15
$ws.Range("O6").Value = 1.67
# Row 8
$ws.Range("G8").Value = 1.45
$ws.Range("I8").Value = 6.5
$ws.Range("N8").Value = 1.7
$ws.Range("O8").Value = 2.1
$ws.Range("R8").Value = 1.83
$ws.Range("S8").Value = 1.83
$ws.Range("T8").Value = 7.5
$ws.Range("W8").Value = 10
$ws.Range("AG8").Value = 19
$ws.Range("AH8").Value = 67
$ws.Range("AI8").Value = 41
$ws.Range("AJ8").Value = 41
# Row 10
$ws.Range("H10").Value = 3.9
$ws.Range("K10").Value = 21
$ws.Range("U10").Value = 26
# Row 11
$ws.Range("G11").Value = 1.48
$ws.Range("R11").Value = 1.73
$ws.Range("S11").Value = 2
$ws.Range("U11").Value = 8
$ws.Range("AA11").Value = 9.5
$ws.Range("AD11").Value = 201
$ws.Range("AE11").Value = 19
$ws.Range("AF11").Value = 34
# Row 12
$ws.Range("G12").Value = 2.7
$ws.Range("I12").Value = 2.4
$ws.Range("L12").Value = 1.17
$ws.Range("M12").Value = 5
$ws.Range("T12").Value = 13
$ws.Range("U12").Value = 17
$ws.Range("V12").Value = 11
$ws.Range("W12").Value = 29
$ws.Range("X12").Value = 19
$ws.Range("Y12").Value = 21
$ws.Range("AC12").Value = 29
$ws.Range("AE12").Value = 13
$ws.Range("AF12").Value = 15
$ws.Range("AG12").Value = 10
$ws.Range("AH12").Value = 23
$ws.Range("AI12").Value = 17
# Row 13
$ws.Range("G13").Value = 2.7
$ws.Range("I13").Value = 2.38
$ws.Range("V13").Value = 11
$ws.Range("W13").Value = 29
$ws.Range("AG13").Value = 9.5
$ws.Range("AJ13").Value = 23
# Row 14
$ws.Range("G14").Value = 1.62
$ws.Range("H14").Value = 4.2
$ws.Range("I14").Value = 4.33
$ws.Range("AE14").Value = 17
$ws.Range("AH14").Value = 51
$ws.Range("AI14").Value = 34
# Row 15
$ws.Range("N15").Value = 1.62
$ws.Range("O15").Value = 2.25
# Row 16
$ws.Range("G16").Value = 1.57
$ws.Range("H16").Value = 4.1
$ws.Range("J16").Value = 1.04
$ws.Range("K16").Value = 13
$ws.Range("N16").Value = 1.75
$ws.Range("O16").Value = 2.05
$ws.Range("R16").Value = 1.8
$ws.Range("S16").Value = 1.91
$ws.Range("T16").Value = 7.5
$ws.Range("X16").Value = 13
$ws.Range("Z16").Value = 13
$ws.Range("AA16").Value = 8
$ws.Range("AD16").Value = 251
$ws.Range("AF16").Value = 26
# Row 19
$ws.Range("G19").Value = 1.91
$ws.Range("H19").Value = 3.75
$ws.Range("I19").Value = 3.3
$ws.Range("N19").Value = 1.62
$ws.Range("O19").Value = 2.25
$ws.Range("U19").Value = 11
$ws.Range("V19").Value = 9
$ws.Range("X19").Value = 15
$ws.Range("AD19").Value = 126
$ws.Range("AE19").Value = 13
$ws.Range("AG19").Value = 12
# Row 20
$ws.Range("G20").Value = 4.5
$ws.Range("I20").Value = 1.8
$ws.Range("J20").Value = 1.07
$ws.Range("K20").Value = 9
$ws.Range("L20").Value = 1.36
$ws.Range("M20").Value = 3
$ws.Range("N20").Value = 2.15
$ws.Range("O20").Value = 1.67
$ws.Range("T20").Value = 11
$ws.Range("W20").Value = 51
$ws.Range("X20").Value = 41
$ws.Range("Z20").Value = 8.5
$ws.Range("AD20").Value = 401
$ws.Range("AG20").Value = 8.5
$ws.Range("AI20").Value = 15
# Row 21
$ws.Range("G21").Value = 2.55
$ws.Range("I21").Value = 2.7
$ws.Range("T21").Value = 9
$ws.Range("AB21").Value = 13
$ws.Range("AD21").Value = 201
$ws.Range("AF21").Value = 15
$ws.Range("AG21").Value = 11
$ws.Range("AH21").Value = 29
# Row 23
$ws.Range("AD23").Value = 1000
# Row 24
$ws.Range("G24").Value = 3.1
$ws.Range("I24").Value = 2.4
$ws.Range("T24").Value = 7.5
$ws.Range("W24").Value = 34
$ws.Range("X24").Value = 29
$ws.Range("Z24").Value = 6.5
$ws.Range("AE24").Value = 6.5
$ws.Range("AG24").Value = 11
$ws.Range("AH24").Value = 23
# Row 25
$ws.Range("N25").Value = 2.5
$ws.Range("O25").Value = 1.5
# Row 28
$ws.Range("G28").Value = 2.55
$ws.Range("I28").Value = 2.75
$ws.Range("K28").Value = 8
$ws.Range("U28").Value = 11
$ws.Range("AE28").Value = 7.5
$ws.Range("AF28").Value = 13
# Row 29
$ws.Range("G29").Value = 1.6
$ws.Range("H29").Value = 3.8
$ws.Range("I29").Value = 4.5
$ws.Range("N29").Value = 1.83
$ws.Range("O29").Value = 2.03
$ws.Range("P29").Value = 1.36
$ws.Range("Q29").Value = 3
$ws.Range("Y29").Value = 26
$ws.Range("Z29").Value = 12
$ws.Range("AA29").Value = 7.5
$ws.Range("AE29").Value = 13
# Row 30
$ws.Range("G30").Value = 3.1
$ws.Range("H30").Value = 3.2
$ws.Range("I30").Value = 2.25
$ws.Range("J30").Value = 1.07
$ws.Range("K30").Value = 6.9
$ws.Range("L30").Value = 1.33
$ws.Range("M30").Value = 3.05
$ws.Range("N30").Value = 1.98
$ws.Range("O30").Value = 1.75
$ws.Range("P30").Value = 1.44
$ws.Range("Q30").Value = 2.62
$ws.Range("R30").Value = 1.75
$ws.Range("S30").Value = 1.95
$ws.Range("V30").Value = 10.75
$ws.Range("X30").Value = 27
$ws.Range("Z30").Value = 6.9
$ws.Range("AA30").Value = 6.2
$ws.Range("AB30").Value = 14
$ws.Range("AC30").Value = 65
$ws.Range("AD30").Value = 500
$ws.Range("AE30").Value = 7.3
$ws.Range("AF30").Value = 10.5
$ws.Range("AG30").Value = 9
$ws.Range("AI30").Value = 19
$ws.Range("AJ30").Value = 30
# Row 31
$ws.Range("G31").Value = 1.78
$ws.Range("H31").Value = 3.65
$ws.Range("I31").Value = 4.05
$ws.Range("L31").Value = 1.23
$ws.Range("M31").Value = 3.75
$ws.Range("N31").Value = 1.7
$ws.Range("O31").Value = 2.05
$ws.Range("U31").Value = 9.5
$ws.Range("V31").Value = 8
$ws.Range("W31").Value = 15
$ws.Range("X31").Value = 13
$ws.Range("AA31").Value = 7.2
$ws.Range("AB31").Value = 13.5
$ws.Range("AC31").Value = 55
$ws.Range("AF31").Value = 24
$ws.Range("AG31").Value = 13
$ws.Range("AH31").Value = 65
$ws.Range("AI31").Value = 35
$ws.Range("AJ31").Value = 37
# Row 32
$ws.Range("G32").Value = 1.72
$ws.Range("H32").Value = 3.95
$ws.Range("I32").Value = 4.1
$ws.Range("R32").Value = 1.7
$ws.Range("S32").Value = 2.05
$ws.Range("AC32").Value = 60
$ws.Range("AE32").Value = 13.5
$ws.Range("AF32").Value = 24
$ws.Range("AH32").Value = 60
# Row 33
$ws.Range("G33").Value = 3.3
$ws.Range("H33").Value = 3.4
$ws.Range("I33").Value = 2.07
$ws.Range("K33").Value = 7.7
$ws.Range("L33").Value = 1.26
$ws.Range("M33").Value = 3.5
$ws.Range("N33").Value = 1.75
$ws.Range("O33").Value = 1.95
$ws.Range("P33").Value = 1.37
$ws.Range("Q33").Value = 2.85
$ws.Range("R33").Value = 1.65
$ws.Range("S33").Value = 2.12
$ws.Range("T33").Value = 11
$ws.Range("U33").Value = 18.5
$ws.Range("V33").Value = 11.25
$ws.Range("W33").Value = 45
$ws.Range("X33").Value = 27
$ws.Range("Z33").Value = 7.7
$ws.Range("AA33").Value = 6.6
$ws.Range("AF33").Value = 11
$ws.Range("AG33").Value = 8.5
$ws.Range("AH33").Value = 20
$ws.Range("AI33").Value = 15.5
$ws.Range("AJ33").Value = 23
# Row 34
$ws.Range("G34").Value = 2.88
$ws.Range("I34").Value = 2.4
$ws.Range("T34").Value = 8.5
$ws.Range("W34").Value = 29
$ws.Range("AE34").Value = 7.5
$ws.Range("AF34").Value = 11
$ws.Range("AH34").Value = 23
$ws.Range("AI34").Value = 21
# Row 36
$ws.Range("N36").Value = 2.4
$ws.Range("O36").Value = 1.53
